$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("G9").Value = 2.45
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 2.8
$ws.Range("J9").Value = 3.25
$ws.Range("L9").Value = 3.6
$ws.Range("M9").Value = 1.05
$ws.Range("O9").Value = 1.37
$ws.Range("W9").Value = 7.5
$ws.Range("AC9").Value = 8
$ws.Range("AE9").Value = 15
$ws.Range("AH9").Value = 13
$ws.Range("AI9").Value = 11
$ws.Range("AJ9").Value = 29
$ws.Range("AK9").Value = 26
$ws.Range("BB9").Value = 201
$ws.Range("G10").Value = 2.4
$ws.Range("H10").Value = 2.7
$ws.Range("I10").Value = 3.5
$ws.Range("L10").Value = 4.33
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 1.57
$ws.Range("P10").Value = 2.25
$ws.Range("Q10").Value = 2.88
$ws.Range("R10").Value = 1.4
$ws.Range("AG10").Value = 7.5
$ws.Range("AI10").Value = 15
$ws.Range("AR10").Value = 101
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 1.3
$ws.Range("P11").Value = 3.4
$ws.Range("Q11").Value = 2.03
$ws.Range("R11").Value = 1.83
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 2.95
$ws.Range("I14").Value = 2.4
$ws.Range("J14").Value = 3.45
$ws.Range("K14").Value = 2
$ws.Range("P14").Value = 2.62
$ws.Range("S14").Value = 1.42
$ws.Range("T14").Value = 2.47
$ws.Range("V14").Value = 1.82
$ws.Range("W14").Value = 8.5
$ws.Range("X14").Value = 15.5
$ws.Range("Z14").Value = 40
$ws.Range("AB14").Value = 35
$ws.Range("AG14").Value = 6.9
$ws.Range("AH14").Value = 11.25
$ws.Range("AJ14").Value = 26
$ws.Range("AL14").Value = 35
$ws.Range("AN14").Value = 4.8
$ws.Range("AP14").Value = 22
$ws.Range("AR14").Value = 100
$ws.Range("AS14").Value = 250
$ws.Range("AT14").Value = 2.42
$ws.Range("AU14").Value = 6.7
$ws.Range("AW14").Value = 4.25
$ws.Range("AX14").Value = 13
$ws.Range("AZ14").Value = 55
$ws.Range("BA14").Value = 100
